$d = $word.ActiveDocument

# Update the date line in the first paragraph
$d.Paragraphs.Item(1).Range.Find.Execute("2025-12-13 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-12-14 Sunday", 1)

# Update each answer cell in the practice table (20 rows x 5 columns)
# Use wdReplaceOne (1) instead of wdReplaceAll so a Find scoped to one cell
# never touches a duplicate value sitting in a different cell.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Find.Execute("11+37=48", $true, $false, $false, $false, $false, $true, 1, $false, "43+56=99", 1)
$t.Cell(1, 2).Range.Find.Execute("91-31=60", $true, $false, $false, $false, $false, $true, 1, $false, "48+35=83", 1)
$t.Cell(1, 3).Range.Find.Execute("56+15=71", $true, $false, $false, $false, $false, $true, 1, $false, "71-62=9", 1)
$t.Cell(1, 4).Range.Find.Execute("28+31=59", $true, $false, $false, $false, $false, $true, 1, $false, "4+57=61", 1)
$t.Cell(1, 5).Range.Find.Execute("44-1=43", $true, $false, $false, $false, $false, $true, 1, $false, "82-24=58", 1)
$t.Cell(2, 1).Range.Find.Execute("94-39=55", $true, $false, $false, $false, $false, $true, 1, $false, "29+15=44", 1)
$t.Cell(2, 2).Range.Find.Execute("55-53=2", $true, $false, $false, $false, $false, $true, 1, $false, "99-31=68", 1)
$t.Cell(2, 3).Range.Find.Execute("71+26=97", $true, $false, $false, $false, $false, $true, 1, $false, "47-46=1", 1)
$t.Cell(2, 4).Range.Find.Execute("2+41=43", $true, $false, $false, $false, $false, $true, 1, $false, "80-11=69", 1)
$t.Cell(2, 5).Range.Find.Execute("69-57=12", $true, $false, $false, $false, $false, $true, 1, $false, "72-22=50", 1)
$t.Cell(3, 1).Range.Find.Execute("55+18=73", $true, $false, $false, $false, $false, $true, 1, $false, "70-4=66", 1)
$t.Cell(3, 2).Range.Find.Execute("30+23=53", $true, $false, $false, $false, $false, $true, 1, $false, "53+36=89", 1)
$t.Cell(3, 3).Range.Find.Execute("44+48=92", $true, $false, $false, $false, $false, $true, 1, $false, "42+14=56", 1)
$t.Cell(3, 4).Range.Find.Execute("44+26=70", $true, $false, $false, $false, $false, $true, 1, $false, "62+16=78", 1)
$t.Cell(3, 5).Range.Find.Execute("27-0=27", $true, $false, $false, $false, $false, $true, 1, $false, "45-4=41", 1)
$t.Cell(4, 1).Range.Find.Execute("49+11=60", $true, $false, $false, $false, $false, $true, 1, $false, "54-8=46", 1)
$t.Cell(4, 2).Range.Find.Execute("19-14=5", $true, $false, $false, $false, $false, $true, 1, $false, "65-18=47", 1)
$t.Cell(4, 3).Range.Find.Execute("93-78=15", $true, $false, $false, $false, $false, $true, 1, $false, "92-19=73", 1)
$t.Cell(4, 4).Range.Find.Execute("60+38=98", $true, $false, $false, $false, $false, $true, 1, $false, "12+78=90", 1)
$t.Cell(4, 5).Range.Find.Execute("88-80=8", $true, $false, $false, $false, $false, $true, 1, $false, "59-52=7", 1)
$t.Cell(5, 1).Range.Find.Execute("24+18=42", $true, $false, $false, $false, $false, $true, 1, $false, "43+48=91", 1)
$t.Cell(5, 2).Range.Find.Execute("92-49=43", $true, $false, $false, $false, $false, $true, 1, $false, "5-4=1", 1)
$t.Cell(5, 3).Range.Find.Execute("74-55=19", $true, $false, $false, $false, $false, $true, 1, $false, "44-22=22", 1)
$t.Cell(5, 4).Range.Find.Execute("68-9=59", $true, $false, $false, $false, $false, $true, 1, $false, "62-39=23", 1)
$t.Cell(5, 5).Range.Find.Execute("20+14=34", $true, $false, $false, $false, $false, $true, 1, $false, "64-46=18", 1)
$t.Cell(6, 1).Range.Find.Execute("90-83=7", $true, $false, $false, $false, $false, $true, 1, $false, "27+17=44", 1)
$t.Cell(6, 2).Range.Find.Execute("99-87=12", $true, $false, $false, $false, $false, $true, 1, $false, "26+20=46", 1)
$t.Cell(6, 3).Range.Find.Execute("57+36=93", $true, $false, $false, $false, $false, $true, 1, $false, "65+5=70", 1)
$t.Cell(6, 4).Range.Find.Execute("56+38=94", $true, $false, $false, $false, $false, $true, 1, $false, "7+68=75", 1)
$t.Cell(6, 5).Range.Find.Execute("90-56=34", $true, $false, $false, $false, $false, $true, 1, $false, "91-70=21", 1)
$t.Cell(7, 1).Range.Find.Execute("83-77=6", $true, $false, $false, $false, $false, $true, 1, $false, "53+3=56", 1)
$t.Cell(7, 2).Range.Find.Execute("41+45=86", $true, $false, $false, $false, $false, $true, 1, $false, "49-18=31", 1)
$t.Cell(7, 3).Range.Find.Execute("51+38=89", $true, $false, $false, $false, $false, $true, 1, $false, "13+4=17", 1)
$t.Cell(7, 4).Range.Find.Execute("50-29=21", $true, $false, $false, $false, $false, $true, 1, $false, "68-10=58", 1)
$t.Cell(7, 5).Range.Find.Execute("78-41=37", $true, $false, $false, $false, $false, $true, 1, $false, "92-77=15", 1)
$t.Cell(8, 1).Range.Find.Execute("65+18=83", $true, $false, $false, $false, $false, $true, 1, $false, "26+48=74", 1)
$t.Cell(8, 2).Range.Find.Execute("30+30=60", $true, $false, $false, $false, $false, $true, 1, $false, "95-54=41", 1)
$t.Cell(8, 3).Range.Find.Execute("8-3=5", $true, $false, $false, $false, $false, $true, 1, $false, "60-6=54", 1)
$t.Cell(8, 4).Range.Find.Execute("56+2=58", $true, $false, $false, $false, $false, $true, 1, $false, "16+15=31", 1)
$t.Cell(8, 5).Range.Find.Execute("54-51=3", $true, $false, $false, $false, $false, $true, 1, $false, "48+1=49", 1)
$t.Cell(9, 1).Range.Find.Execute("90-70=20", $true, $false, $false, $false, $false, $true, 1, $false, "81-71=10", 1)
$t.Cell(9, 2).Range.Find.Execute("19+60=79", $true, $false, $false, $false, $false, $true, 1, $false, "98-94=4", 1)
$t.Cell(9, 3).Range.Find.Execute("89-88=1", $true, $false, $false, $false, $false, $true, 1, $false, "93-13=80", 1)
$t.Cell(9, 4).Range.Find.Execute("66-25=41", $true, $false, $false, $false, $false, $true, 1, $false, "57-3=54", 1)
$t.Cell(9, 5).Range.Find.Execute("26-1=25", $true, $false, $false, $false, $false, $true, 1, $false, "14+4=18", 1)
$t.Cell(10, 1).Range.Find.Execute("59+20=79", $true, $false, $false, $false, $false, $true, 1, $false, "97-18=79", 1)
$t.Cell(10, 2).Range.Find.Execute("0+66=66", $true, $false, $false, $false, $false, $true, 1, $false, "84-78=6", 1)
$t.Cell(10, 3).Range.Find.Execute("34+9=43", $true, $false, $false, $false, $false, $true, 1, $false, "51-2=49", 1)
$t.Cell(10, 4).Range.Find.Execute("75-15=60", $true, $false, $false, $false, $false, $true, 1, $false, "70-13=57", 1)
$t.Cell(10, 5).Range.Find.Execute("22-7=15", $true, $false, $false, $false, $false, $true, 1, $false, "48+16=64", 1)
$t.Cell(11, 1).Range.Find.Execute("56+6=62", $true, $false, $false, $false, $false, $true, 1, $false, "99-8=91", 1)
$t.Cell(11, 2).Range.Find.Execute("75-8=67", $true, $false, $false, $false, $false, $true, 1, $false, "88-38=50", 1)
$t.Cell(11, 3).Range.Find.Execute("89+4=93", $true, $false, $false, $false, $false, $true, 1, $false, "67+22=89", 1)
$t.Cell(11, 4).Range.Find.Execute("1+7=8", $true, $false, $false, $false, $false, $true, 1, $false, "36-32=4", 1)
$t.Cell(11, 5).Range.Find.Execute("27+54=81", $true, $false, $false, $false, $false, $true, 1, $false, "26-18=8", 1)
$t.Cell(12, 1).Range.Find.Execute("52-42=10", $true, $false, $false, $false, $false, $true, 1, $false, "81-56=25", 1)
$t.Cell(12, 2).Range.Find.Execute("41+19=60", $true, $false, $false, $false, $false, $true, 1, $false, "22+20=42", 1)
$t.Cell(12, 3).Range.Find.Execute("60-50=10", $true, $false, $false, $false, $false, $true, 1, $false, "38+32=70", 1)
$t.Cell(12, 4).Range.Find.Execute("21+14=35", $true, $false, $false, $false, $false, $true, 1, $false, "62-36=26", 1)
$t.Cell(12, 5).Range.Find.Execute("70+21=91", $true, $false, $false, $false, $false, $true, 1, $false, "72+7=79", 1)
$t.Cell(13, 1).Range.Find.Execute("38+23=61", $true, $false, $false, $false, $false, $true, 1, $false, "78+13=91", 1)
$t.Cell(13, 2).Range.Find.Execute("37+9=46", $true, $false, $false, $false, $false, $true, 1, $false, "72-10=62", 1)
$t.Cell(13, 3).Range.Find.Execute("79-9=70", $true, $false, $false, $false, $false, $true, 1, $false, "69-35=34", 1)
$t.Cell(13, 4).Range.Find.Execute("73+0=73", $true, $false, $false, $false, $false, $true, 1, $false, "62+33=95", 1)
$t.Cell(13, 5).Range.Find.Execute("3+15=18", $true, $false, $false, $false, $false, $true, 1, $false, "4+67=71", 1)
$t.Cell(14, 1).Range.Find.Execute("74+22=96", $true, $false, $false, $false, $false, $true, 1, $false, "45+10=55", 1)
$t.Cell(14, 2).Range.Find.Execute("97-21=76", $true, $false, $false, $false, $false, $true, 1, $false, "68-24=44", 1)
$t.Cell(14, 3).Range.Find.Execute("41+26=67", $true, $false, $false, $false, $false, $true, 1, $false, "93-53=40", 1)
$t.Cell(14, 4).Range.Find.Execute("19+34=53", $true, $false, $false, $false, $false, $true, 1, $false, "25+67=92", 1)
$t.Cell(14, 5).Range.Find.Execute("99-51=48", $true, $false, $false, $false, $false, $true, 1, $false, "92-68=24", 1)
$t.Cell(15, 1).Range.Find.Execute("30-21=9", $true, $false, $false, $false, $false, $true, 1, $false, "32+15=47", 1)
$t.Cell(15, 2).Range.Find.Execute("87-86=1", $true, $false, $false, $false, $false, $true, 1, $false, "91-22=69", 1)
$t.Cell(15, 3).Range.Find.Execute("66+2=68", $true, $false, $false, $false, $false, $true, 1, $false, "42+32=74", 1)
$t.Cell(15, 4).Range.Find.Execute("37+17=54", $true, $false, $false, $false, $false, $true, 1, $false, "77-6=71", 1)
$t.Cell(15, 5).Range.Find.Execute("87-67=20", $true, $false, $false, $false, $false, $true, 1, $false, "23+36=59", 1)
$t.Cell(16, 1).Range.Find.Execute("14+49=63", $true, $false, $false, $false, $false, $true, 1, $false, "72+16=88", 1)
$t.Cell(16, 2).Range.Find.Execute("52+9=61", $true, $false, $false, $false, $false, $true, 1, $false, "44+13=57", 1)
$t.Cell(16, 3).Range.Find.Execute("91-76=15", $true, $false, $false, $false, $false, $true, 1, $false, "22+12=34", 1)
$t.Cell(16, 4).Range.Find.Execute("40+20=60", $true, $false, $false, $false, $false, $true, 1, $false, "61-26=35", 1)
$t.Cell(16, 5).Range.Find.Execute("95-17=78", $true, $false, $false, $false, $false, $true, 1, $false, "42+40=82", 1)
$t.Cell(17, 1).Range.Find.Execute("54+2=56", $true, $false, $false, $false, $false, $true, 1, $false, "91-46=45", 1)
$t.Cell(17, 2).Range.Find.Execute("21+71=92", $true, $false, $false, $false, $false, $true, 1, $false, "68-45=23", 1)
$t.Cell(17, 3).Range.Find.Execute("33+13=46", $true, $false, $false, $false, $false, $true, 1, $false, "54+33=87", 1)
$t.Cell(17, 4).Range.Find.Execute("7+30=37", $true, $false, $false, $false, $false, $true, 1, $false, "35+11=46", 1)
$t.Cell(17, 5).Range.Find.Execute("45-18=27", $true, $false, $false, $false, $false, $true, 1, $false, "43-36=7", 1)
$t.Cell(18, 1).Range.Find.Execute("75-8=67", $true, $false, $false, $false, $false, $true, 1, $false, "98-8=90", 1)
$t.Cell(18, 2).Range.Find.Execute("85-27=58", $true, $false, $false, $false, $false, $true, 1, $false, "75+6=81", 1)
$t.Cell(18, 3).Range.Find.Execute("13+83=96", $true, $false, $false, $false, $false, $true, 1, $false, "90-38=52", 1)
$t.Cell(18, 4).Range.Find.Execute("63-43=20", $true, $false, $false, $false, $false, $true, 1, $false, "3+96=99", 1)
$t.Cell(18, 5).Range.Find.Execute("6-3=3", $true, $false, $false, $false, $false, $true, 1, $false, "11+59=70", 1)
$t.Cell(19, 1).Range.Find.Execute("14+59=73", $true, $false, $false, $false, $false, $true, 1, $false, "76-57=19", 1)
$t.Cell(19, 2).Range.Find.Execute("29+48=77", $true, $false, $false, $false, $false, $true, 1, $false, "6+69=75", 1)
$t.Cell(19, 3).Range.Find.Execute("46+22=68", $true, $false, $false, $false, $false, $true, 1, $false, "17+70=87", 1)
$t.Cell(19, 4).Range.Find.Execute("18+67=85", $true, $false, $false, $false, $false, $true, 1, $false, "81+14=95", 1)
$t.Cell(19, 5).Range.Find.Execute("99-86=13", $true, $false, $false, $false, $false, $true, 1, $false, "43-23=20", 1)
$t.Cell(20, 1).Range.Find.Execute("17+4=21", $true, $false, $false, $false, $false, $true, 1, $false, "68-30=38", 1)
$t.Cell(20, 2).Range.Find.Execute("50+45=95", $true, $false, $false, $false, $false, $true, 1, $false, "97-94=3", 1)
$t.Cell(20, 3).Range.Find.Execute("6+45=51", $true, $false, $false, $false, $false, $true, 1, $false, "89-35=54", 1)
$t.Cell(20, 4).Range.Find.Execute("54+12=66", $true, $false, $false, $false, $false, $true, 1, $false, "29-25=4", 1)
$t.Cell(20, 5).Range.Find.Execute("30-26=4", $true, $false, $false, $false, $false, $true, 1, $false, "89-33=56", 1)
